$wb = $excel.ActiveWorkbook

# --- Sheet: Pesos_Locales_Económico (sheet11) ---
$ws = $wb.Worksheets.Item("Pesos_Locales_Económico")
$ws.Range("B2").Value = 0.06859393436079969
$ws.Range("B3").Value = 0.06859393436079969
$ws.Range("B4").Value = 0.1451001396860344
$ws.Range("B5").Value = 0.06859393436079965
$ws.Range("B6").Value = 0.1451001396860344
$ws.Range("B7").Value = 0.1451001396860344
$ws.Range("B8").Value = 0.01428744057464852
$ws.Range("B9").Value = 0.01428744057464852
$ws.Range("B10").Value = 0.06859393436079965
$ws.Range("B11").Value = 0.0697866565851307
$ws.Range("B12").Value = 0.008867907233478526
$ws.Range("B13").Value = 0.1451001396860344
$ws.Range("B14").Value = 0.009419377695460395
$ws.Range("B15").Value = 0.01428744057464852
$ws.Range("B16").Value = 0.01428744057464852

# --- Sheet: Ranking_Alternativas (sheet15) ---
$ws = $wb.Worksheets.Item("Ranking_Alternativas")
# Name (alternative) reorder in column A for rows 6,7,8,11,12
$ws.Range("A6").Value = "Quebrada Verde"
$ws.Range("A7").Value = "Laguna Verde"
$ws.Range("A8").Value = "Puertas Negras"
$ws.Range("A11").Value = "Placeres"
$ws.Range("A12").Value = "Reina Isabel 2"
# Updated global-weight values in column B
$ws.Range("B2").Value = 0.1109167152245205
$ws.Range("B3").Value = 0.1003113684231125
$ws.Range("B4").Value = 0.1001482628022656
$ws.Range("B5").Value = 0.0897207308724224
$ws.Range("B6").Value = 0.0895818581737316
$ws.Range("B7").Value = 0.08845323282570873
$ws.Range("B8").Value = 0.08839533147738601
$ws.Range("B9").Value = 0.08768073424173278
$ws.Range("B10").Value = 0.07524167872654772
$ws.Range("B11").Value = 0.03296433373103454
$ws.Range("B12").Value = 0.03255045205159163
$ws.Range("B13").Value = 0.02979481574586582
$ws.Range("B14").Value = 0.02933381661735372
$ws.Range("B15").Value = 0.02555781953583888
$ws.Range("B16").Value = 0.01934884955088798

# --- Sheet: Resultados (sheet2) ---
$ws = $wb.Worksheets.Item("Resultados")
$ws.Range("B2").Value = 0.02933381661735372
$ws.Range("B3").Value = 0.08768073424173278
$ws.Range("B4").Value = 0.02979481574586582
$ws.Range("B5").Value = 0.1001482628022656
$ws.Range("B6").Value = 0.08845323282570873
$ws.Range("B7").Value = 0.0897207308724224
$ws.Range("B8").Value = 0.1003113684231125
$ws.Range("B9").Value = 0.01934884955088798
$ws.Range("B10").Value = 0.03296433373103454
$ws.Range("B11").Value = 0.07524167872654772
$ws.Range("B12").Value = 0.1109167152245205
$ws.Range("B13").Value = 0.08839533147738601
$ws.Range("B14").Value = 0.0895818581737316
$ws.Range("B15").Value = 0.03255045205159163
$ws.Range("B16").Value = 0.02555781953583888

# --- Sheet: Matriz_Económico (sheet6) ---
$ws = $wb.Worksheets.Item("Matriz_Económico")
$ws.Range("D2").Value = 0.3333333333333333
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.3333333333333333
$ws.Range("N2").Value = 7
$ws.Range("D3").Value = 0.3333333333333333
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.3333333333333333
$ws.Range("N3").Value = 7
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 3
$ws.Range("D5").Value = 0.3333333333333333
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.3333333333333333
$ws.Range("N5").Value = 7
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 3
$ws.Range("E6").Value = 3
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 3
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 3
$ws.Range("E7").Value = 3
$ws.Range("J7").Value = 3
$ws.Range("K7").Value = 3
$ws.Range("D10").Value = 0.3333333333333333
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.3333333333333333
$ws.Range("N10").Value = 7
$ws.Range("D11").Value = 0.3333333333333333
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.3333333333333333
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 3
$ws.Range("E13").Value = 3
$ws.Range("J13").Value = 3
$ws.Range("K13").Value = 3
$ws.Range("B14").Value = 0.1428571428571428
$ws.Range("C14").Value = 0.1428571428571428
$ws.Range("E14").Value = 0.1428571428571428
$ws.Range("J14").Value = 0.1428571428571428

Write-Output "Edit complete"